# "Generate Report for Archive"
# - Status text for rows still awaiting translation moves from
#   "Ready for handoff" to "In Translation" on every sheet.
# - The now-shorter status text means the Status columns (Overview!E:F,
#   zh-cn!C, de-de!C) are re-sized narrower to fit the new content.

$wb = $excel.ActiveWorkbook

# Update the status label everywhere it appears (Overview E/F columns,
# and column C on each locale sheet) in one pass per sheet.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# Narrow the Status columns to match the shorter label.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
